# Updates cryptos list values per the scraped diff (prices, % changes,
# plus a rank swap between two coin-pairs) — commit: "Updated cryptos list
# on Mon Aug 19 17:50:54 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few new values are decimals with a significant trailing zero
# ("6.80", "154.80", "0.880"). A plain .Value assignment lets Excel
# auto-convert them to numbers and silently drop the trailing zero,
# so those specific cells are pre-formatted as Text to preserve the
# exact printed string.
foreach ($addr in @("D9", "D32", "D35")) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = [ordered]@{
    'D2' = '58.886.64'
    'D3' = '2.598.57'
    'E3' = '  -1.53%  '
    'D5' = '553.41'
    'E5' = '  +3.23%  '
    'D6' = '144.09'
    'E6' = '  -0.60%  '
    'E7' = '  -0.06%  '
    'D8' = '0.602'
    'E8' = '  +5.46%  '
    'D9' = '6.80'
    'E9' = '  +2.44%  '
    'E10' = '  -1.19%  '
    'D11' = '0.141'
    'E11' = '  +4.48%  '
    'E12' = '  -0.55%  '
    'D13' = '3.054.33'
    'E13' = '  -1.81%  '
    'D14' = '58.822.26'
    'E14' = '  -1.01%  '
    'D15' = '20.87'
    'E15' = '  -1.03%  '
    'D16' = '2.610.06'
    'E16' = '  -0.09%  '
    'E17' = '  -1.71%  '
    'D18' = '4.46'
    'E18' = '  +1.56%  '
    'D19' = '337.39'
    'E19' = '  -0.74%  '
    'D20' = '10.08'
    'E20' = '  -2.62%  '
    'D21' = '6.16'
    'E21' = '  -2.23%  '
    'D22' = '0.999'
    'E22' = '  -0.11%  '
    'D23' = '66.51'
    'E23' = '  -0.86%  '
    'E24' = '  +2.51%  '
    'D25' = '0.996'
    'E25' = '  -0.17%  '
    'D26' = '0.158'
    'E26' = '  -3.85%  '
    'D27' = '7.13'
    'E27' = '  -1.85%  '
    'D28' = '0.0₃0762'
    'E28' = '  +2.29%  '
    'D29' = '0.999'
    'E29' = '  +0.05%  '
    'E30' = '  +1.27%  '
    'D31' = '5.95'
    'E31' = '  +2.09%  '
    'D32' = '154.80'
    'E32' = '  +2.17%  '
    'E33' = '  +0.20%  '
    'D34' = '3.93'
    'E34' = '  -1.64%  '
    'D35' = '0.880'
    'E35' = '  +5.07%  '
    'B36' = 'OKB'
    'C36' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D36' = '37.22'
    'E36' = '  +0.33%  '
    'B37' = 'ImmutableX'
    'C37' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D37' = '1.11'
    'E37' = '  -1.61%  '
    'E38' = '  +0.99%  '
    'E39' = '  -0.89%  '
    'E40' = '  +0.82%  '
    'D41' = '282.32'
    'E41' = '  -2.16%  '
    'D42' = '0.998'
    'E42' = '  -0.17%  '
    'D43' = '0.599'
    'E43' = '  -0.88%  '
    'D44' = '0.0955'
    'E44' = '  +0.97%  '
    'E45' = '  -0.89%  '
    'D46' = '0.0532'
    'E46' = '  -0.30%  '
    'B48' = 'Maker'
    'C48' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D48' = '1.927.04'
    'E48' = '  -2.31%  '
    'B49' = 'RenderToken'
    'C49' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D49' = '4.45'
    'E49' = '  -2.04%  '
    'D50' = '17.89'
    'E50' = '  -2.04%  '
    'D51' = '115.24'
    'E51' = '  +3.70%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Applied $($updates.Count) cell updates."
